$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Resumen de Reunion" meeting-summary docs (15 and 16) for column F
# (Elaboracion Iteracion 2), rows 11-12.
$ws.Range("F11").Value = "Resumen de Reunión 15"
$ws.Range("F12").Value = "Resumen de Reunión 16"

# Column F widened to fit the new text (target stored width 27.85546875
# characters; the host quantizes ColumnWidth writes to 1/6-character
# pixel steps, so 27 is the closest input that round-trips to it).
$ws.Columns.Item(6).ColumnWidth = 27

# Restore the scroll position / selection left by the editing author.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H12").Select()
